# Weekly refresh: a new week's price observation was recorded for
# "Orégano" at Vega Central Mapocho de Santiago, so a new row is inserted
# right after the existing most-recent entry (row 19) and the rest of the
# historical rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 20..41 down to 21..42, freeing up row 20 for the new record.
$ws.Rows("20:20").Insert()

# Populate the newly freed row 20 with the new week's observation.
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44540
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 100112029
$ws.Range("G20").Value = "Orégano"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = 8500
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 8719
$ws.Range("N20").Value = "`$/docena de atados"
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 2906
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = "Hortaliza"
